$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.914.63"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.550.68"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  -0.35%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.38"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  -0.37%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.27"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +3.55%  "
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "1.771.91"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "1.552.71"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "26.904.39"
$ws.Range("E16").Value = "  -0.12%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.59"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "0.0₃0698"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  -0.68%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.94"
$ws.Range("D27").Style = $origStyle
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  +1.72%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.08"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "1.420.24"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("E34").Value = "  +3.75%  "
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +0.49%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.525"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("E43").Value = "  +3.20%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +1.65%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.56"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +1.29%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.76"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "1.685.32"
$ws.Range("E47").Value = "  -0.39%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.29"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("E49").Value = "  +3.95%  "
$ws.Range("E50").Value = "  +1.53%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0960"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +0.52%  "
